# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet and on each per-locale sheet.
#  - Each per-locale sheet gets its "Latest Target File" / "Latest Handback
#    File" / "Latest Handback DateTime" columns populated for both rows, with
#    a new hyperlink on the target-file cell.
#  - A couple of columns are widened to comfortably fit the newly
#    populated long file names.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9b0f12eb571d68ac47545d87bad7b0ad4554548/e2e/0ec1f7f7-388b-4f1c-8058-217002c670f6.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9b0f12eb571d68ac47545d87bad7b0ad4554548/e2e/c87894ca-de31-4dbf-9cc9-3dddda0dc4fe.md"
$mdName1 = "0ec1f7f7-388b-4f1c-8058-217002c670f6.md"
$mdName2 = "c87894ca-de31-4dbf-9cc9-3dddda0dc4fe.md"

# ---- Overview sheet: refresh per-locale status + widen the status columns ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("E1").ColumnWidth = 29.17
$overview.Range("F1").ColumnWidth = 29.17

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, "", "", $mdName1)
$zhcn.Range("J2").Value = "0ec1f7f7-388b-4f1c-8058-217002c670f6.14a395dec7a646555262120416f6376885f570be.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-18 04:46:17"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, "", "", $mdName2)
$zhcn.Range("J3").Value = "c87894ca-de31-4dbf-9cc9-3dddda0dc4fe.ef3495fc278a50e5de3eaecae6afa6d7b2bb0614.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-18 04:46:17"

$zhcn.Range("C1").ColumnWidth = 29.17
$zhcn.Range("I1").ColumnWidth = 39.17
$zhcn.Range("J1").ColumnWidth = 39.17

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, "", "", $mdName1)
$dede.Range("J2").Value = "0ec1f7f7-388b-4f1c-8058-217002c670f6.14a395dec7a646555262120416f6376885f570be.de-de.xlf"
$dede.Range("K2").Value = "2016-10-18 04:46:40"

$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, "", "", $mdName2)
$dede.Range("J3").Value = "c87894ca-de31-4dbf-9cc9-3dddda0dc4fe.ef3495fc278a50e5de3eaecae6afa6d7b2bb0614.de-de.xlf"
$dede.Range("K3").Value = "2016-10-18 04:46:40"

$dede.Range("C1").ColumnWidth = 29.17
$dede.Range("I1").ColumnWidth = 39.17
$dede.Range("J1").ColumnWidth = 39.17
